# Generate Report for Archive
#
# The localization run moved on: the "Ready for handoff" status is now
# "In Translation" for the zh-cn / de-de targets. Update the Status cell
# on every sheet (the per-locale tables plus the rolled-up Overview
# sheet), then right-size the now-narrower Status columns the way Excel
# would after you retype a shorter value into them.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: one column per locale (E = zh-cn, F = de-de) -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns("E:F").ColumnWidth = 12.5

# --- Per-locale detail sheets: Status lives in column C ---------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns("C:C").ColumnWidth = 12.5
